$wb = $excel.ActiveWorkbook

# Sheet 1: "HSV Log OLS"
$ws1 = $wb.Worksheets.Item("HSV Log OLS")
$ws1.Range("A2").Value = 0.118
$ws1.Range("B2").Value = 3.105
$ws1.Range("C2").Value = 0.793

# Sheet 2: "HSV PPML"
$ws2 = $wb.Worksheets.Item("HSV PPML")
$ws2.Range("A2").Value = 0.033
$ws2.Range("B2").Value = 1.149
$ws2.Range("C2").Value = 0.784

# Sheet 3: "HSVT NLLSQ"
$ws3 = $wb.Worksheets.Item("HSVT NLLSQ")
$ws3.Range("A2").Value = -0.092
$ws3.Range("B2").Value = 0.215
$ws3.Range("C2").Value = 14078.66
$ws3.Range("D2").Value = 0.133
